$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("colourmagic")

# Unhide and resize columns A (Shuffle) and B (sentenceID)
$ws1.Columns.Item(1).Hidden = $false
$ws1.Columns.Item(2).Hidden = $false
$ws1.Columns.Item(1).ColumnWidth = 9.166666666666666
$ws1.Columns.Item(2).ColumnWidth = 15.666666666666668

# Make "colourmagic" the active sheet/tab
$ws1.Activate()
